$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 1064 (shifts all following rows down by 2).
$ws.Rows.Item(1064).Insert()
$ws.Rows.Item(1064).Insert()

# New row 1064: Hortaliza Coliflor, Mercado Mayorista Lo Valledor de Santiago, Primera, 05/04/2023 (serial 45021)
$ws.Range("A1064").Value = 6
$ws.Range("B1064").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1064").Value = "Metropolitana"
$ws.Range("D1064").Value2 = 45021
$ws.Range("E1064").Value = 13
$ws.Range("F1064").Value = 100112008
$ws.Range("G1064").Value = "Coliflor"
$ws.Range("H1064").Value = "Sin especificar"
$ws.Range("I1064").Value = "Primera"
$ws.Range("J1064").Value = 11000
$ws.Range("K1064").Value = 800
$ws.Range("L1064").Value = 900
$ws.Range("M1064").Value = 859
$ws.Range("N1064").Value = "$/unidad"
$ws.Range("O1064").Value = "Región Metropolitana"
$ws.Range("P1064").Value = 859
$ws.Range("Q1064").Value = 1
$ws.Range("R1064").Value = "Hortaliza"

# New row 1065: Hortaliza Coliflor, Mercado Mayorista Lo Valledor de Santiago, Segunda, 05/04/2023 (serial 45021)
$ws.Range("A1065").Value = 6
$ws.Range("B1065").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1065").Value = "Metropolitana"
$ws.Range("D1065").Value2 = 45021
$ws.Range("E1065").Value = 13
$ws.Range("F1065").Value = 100112008
$ws.Range("G1065").Value = "Coliflor"
$ws.Range("H1065").Value = "Sin especificar"
$ws.Range("I1065").Value = "Segunda"
$ws.Range("J1065").Value = 5000
$ws.Range("K1065").Value = 500
$ws.Range("L1065").Value = 600
$ws.Range("M1065").Value = 550
$ws.Range("N1065").Value = "$/unidad"
$ws.Range("O1065").Value = "Región Metropolitana"
$ws.Range("P1065").Value = 550
$ws.Range("Q1065").Value = 1
$ws.Range("R1065").Value = "Hortaliza"
